$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# NOTE: cell values below are intentionally written in a specific order so that
# newly introduced shared strings land at the same index as the target workbook
# (new unique strings: Master trapped..., Remote, Local, Force-push..., $ git push -f..., 1. First of all...).

# 1) "Master trapped in a wrong commit tunnel" (new topic header, column B of new row 19)
$ws.Cells.Item(19, 2).Value = "Master trapped in a wrong commit tunnel"

# 2) "Remote" - relabel existing "Git" rows 16 & 17, and the new row 20
$ws.Cells.Item(16, 1).Value = "Remote"
$ws.Cells.Item(17, 1).Value = "Remote"

# 3) "Local" - relabel existing "Git" row 18, and the new row 19
$ws.Cells.Item(18, 1).Value = "Local"
$ws.Cells.Item(19, 1).Value = "Local"

# 4) "Force-push from local to remote" (column B of new row 20)
$ws.Cells.Item(20, 2).Value = "Force-push from local to remote"

# 5) "$ git push -f {target_branch_like_origin} {new_of_local_branch}" (column C of new row 20)
$ws.Cells.Item(20, 3).Value = "`$ git push -f {target_branch_like_origin} {new_of_local_branch}"

# 6) Long explanatory text (column C of new row 19)
$ws.Cells.Item(19, 3).Value = "1. First of all, assume the wrong tunnel commit cannot/shouldnot merge to head, so steps would be delete the master branch first, and then re-create master branch at the right commit`n`$ git checkout master`n`$ git branch wrong_track && git checkout wrong_track  <--- cannot delete master when using it`n`$ git branch -d master`n`$ git checkout {master_or_any_other_commit}`n`$ git branch master"

# 7) Finish row 20 column A label
$ws.Cells.Item(20, 1).Value = "Remote"

# --- Update view selection to the new last-filled cell ---
[void]$ws.Range("C20").Select()

# --- Widen column C slightly to accommodate the new, longer text ---
$ws.Columns.Item(3).ColumnWidth = 83.1
